$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.518.55"
$ws.Range("E2").Value = "'  +1.19%  "
$ws.Range("D3").Value = "'3.145.35"
$ws.Range("E3").Value = "'  +3.74%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'561.43"
$ws.Range("E5").Value = "'  +2.80%  "
$ws.Range("D6").Value = "'145.08"
$ws.Range("E6").Value = "'  +5.90%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("D8").Value = "'3.133.12"
$ws.Range("E8").Value = "'  +3.57%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "'  +2.95%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("E10").Value = "'  +5.10%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("E11").Value = "'  +2.38%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "'  +3.43%  "
$ws.Range("D13").Value = "'37.02"
$ws.Range("E13").Value = "'  +4.71%  "
$ws.Range("D14").Value = "'0.0000223"
$ws.Range("D15").Value = "'3.656.91"
$ws.Range("E15").Value = "'  +3.86%  "
$ws.Range("D16").Value = "'64.564.63"
$ws.Range("E16").Value = "'  +1.28%  "
$ws.Range("D17").Value = "'3.162.40"
$ws.Range("E17").Value = "'  +3.99%  "
$ws.Range("E18").Value = "'  +1.82%  "
$ws.Range("D19").Value = "'514.27"
$ws.Range("E19").Value = "'  +6.90%  "
$ws.Range("D20").Value = "'6.87"
$ws.Range("E20").Value = "'  +5.49%  "
$ws.Range("D21").Value = "'14.08"
$ws.Range("E21").Value = "'  +3.89%  "
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "'  +6.14%  "
$ws.Range("D23").Value = "'7.47"
$ws.Range("E23").Value = "'  +6.06%  "
$ws.Range("D24").Value = "'13.00"
$ws.Range("E24").Value = "'  +5.37%  "
$ws.Range("D25").Value = "'78.68"
$ws.Range("E25").Value = "'  +1.32%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "'  +19.13%  "
$ws.Range("D28").Value = "'2.84"
$ws.Range("E28").Value = "'  +5.78%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "'  +6.02%  "
$ws.Range("E30").Value = "'  +0.05%  "
$ws.Range("D31").Value = "'26.53"
$ws.Range("E31").Value = "'  +3.97%  "
$ws.Range("D32").Value = "'2.62"
$ws.Range("E32").Value = "'  +0.44%  "
$ws.Range("D33").Value = "'1.14"
$ws.Range("E33").Value = "'  +3.79%  "
$ws.Range("D34").Value = "'545.78"
$ws.Range("E34").Value = "'  -4.32%  "
$ws.Range("D35").Value = "'5.42"
$ws.Range("E35").Value = "'  +2.20%  "
$ws.Range("D36").Value = "'6.08"
$ws.Range("E36").Value = "'  +4.95%  "
$ws.Range("D37").Value = "'54.06"
$ws.Range("E37").Value = "'  +4.91%  "
$ws.Range("D38").Value = "'0.0435"
$ws.Range("E38").Value = "'  +6.89%  "
$ws.Range("D39").Value = "'0.0829"
$ws.Range("E39").Value = "'  +5.47%  "
$ws.Range("D40").Value = "'3.147.33"
$ws.Range("E40").Value = "'  +7.84%  "
$ws.Range("E41").Value = "'  +6.09%  "
$ws.Range("D42").Value = "'2.78"
$ws.Range("E42").Value = "'  -0.16%  "
$ws.Range("D43").Value = "'8.32"
$ws.Range("E43").Value = "'  +2.15%  "
$ws.Range("D44").Value = "'0.270"
$ws.Range("E44").Value = "'  +12.76%  "
$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = "'  +8.98%  "
$ws.Range("D47").Value = "'25.66"
$ws.Range("E47").Value = "'  +4.78%  "
$ws.Range("D48").Value = "'122.48"
$ws.Range("E48").Value = "'  +4.16%  "
$ws.Range("D49").Value = "'0.0₃0525"
$ws.Range("E49").Value = "'  +0.03%  "
$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "'  +1.51%  "
$ws.Range("D51").Value = "'2.12"
$ws.Range("E51").Value = "'  +4.70%  "
